# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the per-locale
# handback-status report sheets ("zh-cn" and "de-de").

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 03:00:53"
$wsZhCn.Range("H2").Value = "2016-03-24 03:02:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 03:01:01"
$wsDeDe.Range("H2").Value = "2016-03-24 03:02:24"
